$d = $word.ActiveDocument

# New paragraphs appended after the final existing paragraph
# ("Talk about how it's given me the ability to drive myself into any
# direction of interest, which makes me...")
$newParagraphs = @(
    "The days when I just can’t",
    "Have you ever woken up and wished that you hadn’t? Not that you had more time to sleep, but legitimately wishing that you could continue sleeping forever. Not having to worry about your schedule, what you’re planning to get done during the week, what you need to eat, to water your plants, etc. Nearing the end of my PhD, my mind is tired and restless. I find myself able to think at a high level in smaller spurts than in the past, and I am unable to push through those boundaries of my mind as easily. I feel dissociated from my being, pounded down by the weight of the expectations. The closer the finish line gets, the more weight: Get that one last piece of data that will make your story sell itself, that will make it make even more sense.",
    "What are you supposed to do on the days where you just don’t feel like yourself?",
    "The days where you can’t get up out of bed. Those days when hearing yourself think in your head leaves you with a throbbing headache. Where the only thing that allows your mind to numb is the relaxation found within sleep, where there are no expectations but to do rest. ",
    "I want to rest. Like actually rest. To not feel like I need to be up working, ",
    "There are some days when I just can’t"
)

foreach ($text in $newParagraphs) {
    $lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
    $r = $lastPara.Range
    $r.Collapse(0)
    $r.InsertParagraphAfter()

    $newLastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
    $nr = $newLastPara.Range
    $nr.Collapse(0)
    $nr.InsertAfter($text)
}
